# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update @ 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitUpdates = @{
    6  = 2881
    8  = 1897
    11 = 775
    12 = 921
    19 = 6876
    21 = 1659
    22 = 169
    25 = 326
    28 = 1112
    33 = 800
    34 = 1921
    37 = 235
    38 = 28
    39 = 146
    40 = 238
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allUpdates = @{
    9  = 2881
    11 = 1897
    14 = 775
    16 = 921
    22 = 6876
    24 = 1659
    26 = 169
    29 = 326
    32 = 1112
    37 = 800
    38 = 1921
    41 = 235
    42 = 28
    43 = 146
    44 = 238
    49 = 174
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

$wb.Save()
